$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.924.85'
$ws.Range('E2').Value = '  -1.17%  '

$ws.Range('D3').Value = '3.425.13'
$ws.Range('E3').Value = '  +3.45%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '256.47'
$ws.Range('E5').Value = '  +0.42%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '656.42'
$ws.Range('E6').Value = '  +5.18%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.47'
$ws.Range('E7').Value = '  +0.36%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.429'
$ws.Range('E8').Value = '  +5.13%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.05'
$ws.Range('E9').Value = '  +8.10%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.998'
$ws.Range('E10').Value = '  -0.10%  '

$ws.Range('D11').Value = '3.424.13'
$ws.Range('E11').Value = '  +3.54%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.214'
$ws.Range('E12').Value = '  +6.37%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '42.44'
$ws.Range('E13').Value = '  +6.85%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.55'
$ws.Range('E14').Value = '  +19.37%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000259'
$ws.Range('E15').Value = '  +3.50%  '

$ws.Range('D16').Value = '97.696.85'
$ws.Range('E16').Value = '  -1.21%  '

$ws.Range('D17').Value = '4.056.02'
$ws.Range('E17').Value = '  +3.14%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.66'
$ws.Range('E18').Value = '  +36.97%  '

$ws.Range('D19').Value = '3.405.33'
$ws.Range('E19').Value = '  +2.53%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.63'
$ws.Range('E20').Value = '  +13.24%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.520'
$ws.Range('E21').Value = '  +68.16%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.94'
$ws.Range('E22').Value = '  +15.70%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.47'
$ws.Range('E23').Value = '  -0.17%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '510.22'
$ws.Range('E24').Value = '  +4.80%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000206'
$ws.Range('E25').Value = '  +1.52%  '

$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '99.00'
$ws.Range('E26').Value = '  +11.17%  '

$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.09'
$ws.Range('E27').Value = '  +7.99%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.73'
$ws.Range('E28').Value = '  +6.16%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.152'
$ws.Range('E29').Value = '  +11.60%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '11.40'
$ws.Range('E30').Value = '  +10.77%  '

$ws.Range('E31').Value = '  -0.01%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.196'
$ws.Range('E32').Value = '  +4.46%  '

$ws.Range('E33').Value = '  +0.13%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.575'
$ws.Range('E34').Value = '  +21.44%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '29.85'
$ws.Range('E35').Value = '  +7.07%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.20'
$ws.Range('E36').Value = '  +12.97%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '7.86'
$ws.Range('E37').Value = '  +8.81%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.157'
$ws.Range('E38').Value = '  +6.12%  '

$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.40'
$ws.Range('E39').Value = '  +13.64%  '

$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '516.08'
$ws.Range('E40').Value = '  +5.12%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '24.71'
$ws.Range('E41').Value = '  -0.50%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.856'
$ws.Range('E42').Value = '  +8.74%  '

$ws.Range('E43').Value = '  +26.18%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.67'
$ws.Range('E44').Value = '  +0.75%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.31'
$ws.Range('E45').Value = '  +6.08%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.40'
$ws.Range('E46').Value = '  +14.31%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.18'
$ws.Range('E48').Value = '  +11.61%  '

$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.59'
$ws.Range('E49').Value = '  +16.81%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.10'
$ws.Range('E50').Value = '  +7.74%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '50.86'
$ws.Range('E51').Value = '  +10.67%  '
